# Weekly update to the "Hortaliza, Terminal Hortofrutícola Agro Chillán - Pimiento"
# sheet: a new price-report row is inserted at row 97 (pushing every
# subsequent row down by one, so the former last row 197 becomes row 198).
#
# Excel's own "insert entire row" behaviour copies the formatting of the
# row above down into the freshly inserted row (this is what keeps column D's
# date-number-format style intact for the new row), so we let
# Rows.Item(97).Insert() do that instead of touching styles by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 97:197 down to 98:198, leaving a blank row 97 to fill in.
$ws.Rows.Item(97).Insert()

$ws.Cells.Item(97, 1).Value  = 7
$ws.Cells.Item(97, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(97, 3).Value  = "Ñuble"
$ws.Cells.Item(97, 4).Value  = 44546
$ws.Cells.Item(97, 5).Value  = 16
$ws.Cells.Item(97, 6).Value  = 100112002
$ws.Cells.Item(97, 7).Value  = "Pimiento"
$ws.Cells.Item(97, 8).Value  = "Cuatro cascos verde"
$ws.Cells.Item(97, 9).Value  = "Primera"
$ws.Cells.Item(97, 10).Value = 240
$ws.Cells.Item(97, 11).Value = 13000
$ws.Cells.Item(97, 12).Value = 13500
$ws.Cells.Item(97, 13).Value = 13250
$ws.Cells.Item(97, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(97, 15).Value = "Región del Maule"
$ws.Cells.Item(97, 16).Value = 883
$ws.Cells.Item(97, 17).Value = 15
$ws.Cells.Item(97, 18).Value = "Hortaliza"
